{"js": "// The \"competitions\" schema table (Relational Models section) lists the\n// column types for `start_time` / `end_time` as \"datetime\"; the edit\n// shortens both of those Type-column entries to \"date\" (the paragraph\n// marks / paragraph formatting stay untouched \u2014 only the run text\n// changes, collapsing each paragraph down to a single run).\nconst results = context.document.body.search(\"datetime\", {\n  matchCase: true,\n  matchWholeWord: true\n});\nresults.load(\"items\");\nawait context.sync();\n\nfor (const result of results.items) {\n  result.insertText(\"date\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The \"competitions\" schema table (Relational Models section) lists the\n# column types for `start_time` / `end_time` as \"datetime\"; shorten both\n# of those Type-column entries to \"date\" (paragraph marks / paragraph\n# formatting are left untouched \u2014 only the run text changes).\n$d = $word.ActiveDocument\n\n$result = $d.Content.Find.Execute(\n    \"datetime\",   # FindText\n    $false,       # MatchCase\n    $true,        # MatchWholeWord\n    $false,       # MatchWildcards\n    $false,       # MatchSoundsLike\n    $false,       # MatchAllWordForms\n    $true,        # Forward\n    1,            # Wrap (wdFindContinue)\n    $false,       # Format\n    \"date\",       # ReplaceWith\n    2             # Replace (wdReplaceAll)\n)\n"}
